$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell text content (rows reordered/corrected, row 13 removed) ---

$ws.Range("B2").Value = "The COVID-19 vaccine rollout begins with phase 1A which included  healthcare workers, communities in long term care facilities and intellectual disabilities care facilities. "
$ws.Rows.Item(2).RowHeight = 17

$ws.Range("B3").Value = "mRNA-1273 (Moderna) receive Emergency Use Authorization (EUA)."
$ws.Rows.Item(3).RowHeight = 17

$ws.Range("B4").Value = "Phase 1B commences with the vaccination of adults 65 and  older. "
$ws.Rows.Item(4).RowHeight = 17

$ws.Range("B5").Value = "Puerto Rico's Secretary of Health signs administrative order establishing that for the next 28 days, vaccination will be  exclusive for those 65 and older. "
$ws.Rows.Item(5).RowHeight = 17

$ws.Range("B6").Value = "Ad26.COV2.S (Johnson & Johnson) recieves EUA"
$ws.Rows.Item(6).RowHeight = 17

$ws.Range("B7").Value = "On March 3rd, 2021 Secretary of Health  signs  administrative order establishing that starting on March 11th  and for the following 30 days, first doses are to be administered exclusively to adults 60 and older with certain chronic conditions. "
$ws.Rows.Item(7).RowHeight = 34

$ws.Range("B8").Value = "On March 10th, 2021,  Secretary signs  administrative order establishing that starting on March 11th  and for the following 30 days, first doses are to be administered  exclusively to adults 60 and older and 50 to 59 year olds  with  chronic conditions. "
$ws.Rows.Item(8).RowHeight = 34

$ws.Range("B9").Value = "Secretary of Health signs executive order authorizing the  vaccination of personnel in food industry, drug companies, medical  equipment, the public transport sector, air transport and   maritime cargo."
$ws.Rows.Item(9).RowHeight = 17

$ws.Range("B10").Value = "Phase 1C begins with the vaccination of people 50 and older  and 35 and older with chronic conditions. "
$ws.Rows.Item(10).RowHeight = 17

$ws.Range("B11").Value = "Phase 2 begins with vaccination available to everyone 16 and older. "
$ws.Rows.Item(11).RowHeight = 17

# Row 12: date changes from 44302 to 44328, text becomes CDC recommends... (formerly row13 content)
$ws.Range("A12").Value = 44328
$ws.Range("B12").Value = "CDC recommends vaccination for people 12 years and older. "
$ws.Rows.Item(12).RowHeight = 17

# Remove old row 13 (now redundant / duplicate after content shift)
$ws.Rows.Item(13).Delete()

# --- Column width change for column B ---
# (engine rounds stored width to ColumnWidth + 5/6, so back the input off by
# that fixed offset to land exactly on the target stored width of 172.5)
$ws.Columns.Item(2).ColumnWidth = 171.66666666666666

# --- Selection / view changes ---
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("A12:B12").Select()
